$d = $word.ActiveDocument
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Insert-ParaBefore($paraIndex, $xmlInner) {
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $p.Range
    $ins = $d.Range($rng.Start, $rng.Start)
    $xml = "<w:p $ns>$xmlInner</w:p>"
    $ins.InsertXML($xml)
}

function Delete-Para($paraIndex) {
    $p = $d.Paragraphs.Item($paraIndex)
    $p.Range.Delete()
}

# Paragraphs are processed from the bottom of the document upward so that
# earlier (lower-numbered) paragraph indices stay valid while later ones
# are being inserted/removed.

# ---- Change 4 (UC06 "Pagar pedido"): "Cliente realiza o P" + "IX" -> single run ----
$pPr4 = '<w:pPr><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:spacing w:line="240" w:lineRule="auto"/></w:pPr>'
$body4 = '<w:r><w:t>Cliente realiza o PIX</w:t></w:r>'
Insert-ParaBefore 172 "$pPr4$body4"
Delete-Para 173

# ---- Change 3 (UC05 "Visualizar o pedido", Fluxo Secundario): merge split run ----
$pPr3 = '<w:pPr><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr><w:spacing w:line="240" w:lineRule="auto"/></w:pPr>'
$body3 = '<w:r><w:t>A lanchonete vai abrir o aplicativo e selecionar a op' + [char]0x00E7 + [char]0x00E3 + 'o pedidos.</w:t></w:r>'
Insert-ParaBefore 145 "$pPr3$body3"
Delete-Para 146

# ---- Change 2 (UC03 "Agendar pedido", Fluxo Principal, 2nd bullet): reword +
# split into runs with a spellStart/spellEnd proofErr pair, and add a new
# bullet paragraph right after it ----
$pPr2 = '<w:pPr><w:pStyle w:val="PargrafodaLista"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:spacing w:line="240" w:lineRule="auto"/></w:pPr>'
$body2a = '<w:r><w:t xml:space="preserve">O aplicativo </w:t></w:r><w:r><w:t xml:space="preserve">notifica </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>o</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> cliente que o pedido foi realizado com sucesso</w:t></w:r><w:r><w:t>.</w:t></w:r>'
$body2b = '<w:r><w:t>O aplicativo notifica a lanchonete que um novo pedido foi realizado.</w:t></w:r>'
Insert-ParaBefore 82 "$pPr2$body2a"
Insert-ParaBefore 83 "$pPr2$body2b"
Delete-Para 84

# ---- Change 1 (UC03 "Agendar pedido", Fluxo Principal, 1st bullet): split
# into runs with a gramStart/gramEnd proofErr pair ----
$pPr1 = '<w:pPr><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:spacing w:line="240" w:lineRule="auto"/></w:pPr>'
$body1 = '<w:r><w:t xml:space="preserve">O cliente seleciona os itens do </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>pedido .</w:t></w:r><w:proofErr w:type="gramEnd"/>'
Insert-ParaBefore 81 "$pPr1$body1"
Delete-Para 82

Write-Output "Done."
